$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The underlying (player, position, team) triples are unchanged; only the
# row order of the roster table (A2:C19) was shuffled by the author when
# the file was re-uploaded. Re-write the data block with the rows in the
# new order.
$data = @(
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Andrew Wiggins", "SF,PF", "Golden State Warriors"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Malik Monk", "SG,SF", "Sacramento Kings"),
    @("Wendell Carter Jr.", "C", "Orlando Magic"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Anfernee Simons", "PG,SG", "Portland Trail Blazers"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Draymond Green", "PF,C", "Golden State Warriors"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers"),
    @("Aaron Gordon", "PF,C", "Denver Nuggets"),
    @("Kyle Kuzma", "PF", "Washington Wizards")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
